$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = -11.8084
$ws.Range("B7").Value = 4.690799999999999
$ws.Range("E7").Value = 16.11180000000001
$ws.Range("A9").Value = -21.6277
$ws.Range("E10").Value = 16.79119999999999
$ws.Range("B12").Value = 5.411799999999999
$ws.Range("E13").Value = 16.58530000000001
$ws.Range("B14").Value = 6.064699999999997
$ws.Range("C15").Value = -13.67149999999999
$ws.Range("E16").Value = 15.95900000000001
$ws.Range("A18").Value = -22.24470000000002
$ws.Range("A20").Value = -19.50179999999998
$ws.Range("E20").Value = 16.15879999999999
$ws.Range("E24").Value = 16.5736
$ws.Range("B26").Value = 4.074800000000006
$ws.Range("A27").Value = -21.55719999999997
$ws.Range("B27").Value = 5.385800000000004
$ws.Range("B29").Value = 4.824599999999998
$ws.Range("C33").Value = -11.29179999999999
$ws.Range("A35").Value = -19.24669999999999
$ws.Range("C35").Value = -11.5487
$ws.Range("B37").Value = 9.074000000000005
$ws.Range("B38").Value = 4.462000000000002
$ws.Range("C38").Value = -12.0033
$ws.Range("E39").Value = 16.0696
$ws.Range("C43").Value = -13.28589999999999
$ws.Range("C44").Value = -13.34899999999999
$ws.Range("C47").Value = -12.17689999999999
$ws.Range("E47").Value = 16.82929999999999
$ws.Range("E48").Value = 17.5246
$ws.Range("B51").Value = 6.219800000000005
$ws.Range("C51").Value = -11.5547
$ws.Range("B52").Value = 5.085999999999999
$ws.Range("E52").Value = 17.05960000000001
$ws.Range("B55").Value = 4.859399999999997
$ws.Range("E56").Value = 16.58740000000001
$ws.Range("C57").Value = -13.61759999999999
$ws.Range("C63").Value = -12.44090000000001
$ws.Range("A69").Value = -21.6355
$ws.Range("B69").Value = 5.235999999999999
$ws.Range("B70").Value = 6.292900000000007
$ws.Range("C70").Value = -11.5211
$ws.Range("A76").Value = -20.00879999999998
$ws.Range("A78").Value = -20.15879999999998
$ws.Range("B81").Value = 5.4309
$ws.Range("A82").Value = -21.8168
$ws.Range("A83").Value = -22.0405
$ws.Range("B83").Value = 6.2474
$ws.Range("E84").Value = 16.79799999999999
$ws.Range("C88").Value = -11.9643
$ws.Range("A93").Value = -20.46009999999998
$ws.Range("C99").Value = -12.19429999999999
$ws.Range("E100").Value = 16.47230000000001
$ws.Range("E101").Value = 16.81220000000001
$ws.Range("B102").Value = 9.062000000000006
